$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.084.98"
Set-TextValue $ws.Range("E2") "  -0.74%  "
Set-TextValue $ws.Range("D3") "1.668.80"
Set-TextValue $ws.Range("E3") "  -0.34%  "
Set-TextValue $ws.Range("E4") "  -0.18%  "
Set-TextValue $ws.Range("D5") "216.11"
Set-TextValue $ws.Range("E5") "  -0.99%  "
Set-TextValue $ws.Range("D6") "0.5105"
Set-TextValue $ws.Range("E6") "  +0.27%  "
Set-TextValue $ws.Range("E7") "  -0.22%  "
Set-TextValue $ws.Range("E8") "  +1.36%  "
Set-TextValue $ws.Range("D9") "0.06375"
Set-TextValue $ws.Range("E9") "  +0.89%  "
Set-TextValue $ws.Range("D10") "21.83"
Set-TextValue $ws.Range("E10") "  -0.92%  "
Set-TextValue $ws.Range("D11") "0.07439"
Set-TextValue $ws.Range("E11") "  +1.13%  "
Set-TextValue $ws.Range("D12") "1.668.51"
Set-TextValue $ws.Range("E12") "  -0.46%  "
Set-TextValue $ws.Range("D13") "4.508"
Set-TextValue $ws.Range("E13") "  -0.59%  "
Set-TextValue $ws.Range("D14") "0.5811"
Set-TextValue $ws.Range("E14") "  +1.36%  "
Set-TextValue $ws.Range("D15") "0.000008494"
Set-TextValue $ws.Range("E15") "  -0.39%  "
Set-TextValue $ws.Range("D16") "64.10"
Set-TextValue $ws.Range("E16") "  -0.86%  "
Set-TextValue $ws.Range("D17") "25.893.73"
Set-TextValue $ws.Range("E17") "  -1.84%  "
Set-TextValue $ws.Range("D18") "4.928"
Set-TextValue $ws.Range("E18") "  -1.32%  "
Set-TextValue $ws.Range("E19") "  -0.11%  "
Set-TextValue $ws.Range("D20") "10.80"
Set-TextValue $ws.Range("E20") "  -0.33%  "
Set-TextValue $ws.Range("D21") "189.40"
Set-TextValue $ws.Range("E21") "  +1.97%  "
Set-TextValue $ws.Range("D22") "6.187"
Set-TextValue $ws.Range("E22") "  -0.45%  "
Set-TextValue $ws.Range("E23") "  -0.26%  "
Set-TextValue $ws.Range("D24") "144.71"
Set-TextValue $ws.Range("E24") "  +0.86%  "
Set-TextValue $ws.Range("D25") "7.592"
Set-TextValue $ws.Range("E25") "  +0.94%  "
Set-TextValue $ws.Range("D26") "0.1225"
Set-TextValue $ws.Range("E26") "  +4.63%  "
Set-TextValue $ws.Range("D27") "15.67"
Set-TextValue $ws.Range("E27") "  -0.20%  "
Set-TextValue $ws.Range("D28") "0.06658"
Set-TextValue $ws.Range("E28") "  +14.55%  "
Set-TextValue $ws.Range("D29") "1.339"
Set-TextValue $ws.Range("E29") "  +0.56%  "
Set-TextValue $ws.Range("E30") "  -0.95%  "
Set-TextValue $ws.Range("D31") "3.571"
Set-TextValue $ws.Range("E31") "  +1.91%  "
Set-TextValue $ws.Range("D32") "3.523"
Set-TextValue $ws.Range("E32") "  +0.71%  "
Set-TextValue $ws.Range("D33") "1.662"
Set-TextValue $ws.Range("E33") "  +0.88%  "
Set-TextValue $ws.Range("E34") "  +1.23%  "
Set-TextValue $ws.Range("D35") "0.6165"
Set-TextValue $ws.Range("E35") "  +3.85%  "
Set-TextValue $ws.Range("E36") "  +0.41%  "
Set-TextValue $ws.Range("D37") "2.686"
Set-TextValue $ws.Range("E37") "  +0.64%  "
Set-TextValue $ws.Range("D38") "6.248"
Set-TextValue $ws.Range("E38") "  +6.08%  "
Set-TextValue $ws.Range("D39") "1.094.09"
Set-TextValue $ws.Range("E39") "  -0.34%  "
Set-TextValue $ws.Range("D40") "0.01600"
Set-TextValue $ws.Range("E40") "  -0.17%  "
Set-TextValue $ws.Range("D41") "0.8708"
Set-TextValue $ws.Range("E41") "  +1.18%  "
Set-TextValue $ws.Range("E42") "  +0.78%  "
Set-TextValue $ws.Range("D43") "101.21"
Set-TextValue $ws.Range("E43") "  +1.52%  "
Set-TextValue $ws.Range("D44") "1.815.28"
Set-TextValue $ws.Range("E44") "  -1.08%  "
Set-TextValue $ws.Range("D45") "0.00000000116"
Set-TextValue $ws.Range("E45") "  +1.82%  "
Set-TextValue $ws.Range("E46") "  +0.38%  "
Set-TextValue $ws.Range("B47") "EnergySwap"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D47") "8.133"
Set-TextValue $ws.Range("E47") "  +1.17%  "
Set-TextValue $ws.Range("B48") "Frax"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D48") "1.003"
Set-TextValue $ws.Range("E48") "  -0.03%  "
Set-TextValue $ws.Range("D49") "0.05233"
Set-TextValue $ws.Range("E49") "  +0.53%  "
Set-TextValue $ws.Range("D50") "0.4279"
Set-TextValue $ws.Range("E50") "  -0.82%  "
Set-TextValue $ws.Range("D51") "5.988"
Set-TextValue $ws.Range("E51") "  +2.71%  "
